# "add live and ani wood brick"
#
# The document has an empty paragraph (style "Normal", 360-twip / 18pt left
# indent) that sits right before the "Các trạng thái của goomba:" Heading 2.
# That empty paragraph becomes a new "Type:" Heading 2 (numbered list,
# numId 9) followed by a new ListParagraph bullet (numId 1) describing the
# PARA_GOOMBA_BROWN constant.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "`r" -and $p.Style.NameLocal -eq "Normal" -and $p.LeftIndent -eq 18) {
        $next = $p.Next()
        if ($next -ne $null -and $next.Range.Text -match "trạng thái của goomba") {
            $target = $p
            break
        }
    }
}

if ($target -eq $null) {
    throw "Could not locate the empty placeholder paragraph before 'Các trạng thái của goomba:'"
}

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr></w:pPr><w:r><w:t>Type:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="6F008A"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>PARA_GOOMBA_BROWN</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="6F008A"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve"> = 2: màu nâu</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.Range.InsertXML($xml) | Out-Null

Write-Host "Replaced placeholder paragraph with 'Type:' heading and PARA_GOOMBA_BROWN entry."
